$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update translations (column C = en-US) ---
# Row 2: key 10013440 - "interact" -> "demo"
$ws.Range("C2").Value = "demo"

# The following rows' English placeholder translations were removed
# (cleared back to blank, matching the "pending" status rows).
$ws.Range("C3").Value = ""    # 10141355 "switching..."
$ws.Range("C19").Value = ""   # 15929414 "joinned time"
$ws.Range("C20").Value = ""   # 21511119 "curretn lang"
$ws.Range("C24").Value = ""   # 41215122 "message notify"
$ws.Range("C27").Value = ""   # 49032491 "publish article"
$ws.Range("C28").Value = ""   # 52134415 "highly user"
$ws.Range("C29").Value = ""   # 61471412 "edit info"
$ws.Range("C30").Value = ""   # 71101173 "database display"

# --- Update the active selection shown when the sheet is reopened ---
$ws.Range("F10").Select() | Out-Null
